$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new data rows (261-263) to the existing "date / weekday / hour /
# ranking" table, continuing on from the last existing row (260).
$newRows = @(
    @("2025/11/15", "土", 18, 201),
    @("2025/11/15", "土", 19, 201),
    @("2025/11/15", "土", 20, 201)
)

$startRow = 261
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A holds a date-like string ("YYYY/MM/DD"); force it to be stored
    # as text (matching the rest of the column) instead of letting Excel
    # auto-convert it to a date serial number. Resetting to the "Normal"
    # style afterwards avoids leaving a stray number-format override behind.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $rowData[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
